$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New localization rows appended to the bottom of the "Compose Mail" table,
# adding strings for the blank-history / history-record features.
# Written in the same order the author entered them so the shared-string
# table gets appended in matching order: A101, B101, A102, B102, C102, C101.

$ws.Range("A101").Value = "lang_blank_history"
$ws.Range("B101").Value = "Lịch Sử Trống"

$ws.Range("A102").Value = "lang_record"
# Leading apostrophe forces Excel's text-quote-prefix (cell starts with "_"
# so it is entered as literal text), matching the quotePrefix style used
# in the target workbook.
$ws.Range("B102").Value = "'_from_ đến _to_"

$ws.Range("C102").Value = "_from_  to  _to_"

$ws.Range("C101").Value = "History Is Empty"

# Move the view to show the newly added rows, matching the saved selection.
$ws.Range("C101").Select()
